# Update "想去人数" (interest count) figures in the 广州-漫展信息 workbook.
# Sheet order: 1=展览 (Exhibition), 2=演出 (Performance), 3=本地生活 (Local life), 4=全部类型 (All types)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 664
$ws1.Range("F4").Value  = 1311
$ws1.Range("F6").Value  = 71
$ws1.Range("F11").Value = 957
$ws1.Range("F12").Value = 303
$ws1.Range("F13").Value = 163
$ws1.Range("F16").Value = 356
$ws1.Range("F17").Value = 327
$ws1.Range("F18").Value = 741
$ws1.Range("F19").Value = 118
$ws1.Range("F20").Value = 696
$ws1.Range("F21").Value = 240
$ws1.Range("F24").Value = 418
$ws1.Range("F27").Value = 342
$ws1.Range("F30").Value = 446

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 348
$ws2.Range("F5").Value = 30

# --- Sheet 4: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value  = 664
$ws4.Range("F5").Value  = 1311
$ws4.Range("F8").Value  = 71
$ws4.Range("F13").Value = 957
$ws4.Range("F14").Value = 303
$ws4.Range("F15").Value = 163
$ws4.Range("F18").Value = 348
$ws4.Range("F20").Value = 30
$ws4.Range("F21").Value = 356
$ws4.Range("F24").Value = 327
$ws4.Range("F25").Value = 741
$ws4.Range("F26").Value = 118
$ws4.Range("F27").Value = 696
$ws4.Range("F28").Value = 240
$ws4.Range("F31").Value = 418
$ws4.Range("F36").Value = 342
$ws4.Range("F42").Value = 446

Write-Output "Updated 34 cells across sheets 1, 2, and 4."
